# Category matching split as list
# Adds four new per-category sheets ("23 - 9 ...") mirroring the header-only
# layout used for every other "<week> - 9 <brand>" sheet in this workbook.

$wb = $excel.ActiveWorkbook

# Every tracker sheet shares the same 4-column header row
# (Hora, Cambió, Nuevos, Actualizados) styled with the bold/bordered header
# style. Use the most recent existing sheet as the header template.
$template = $wb.Worksheets.Item("22 - 9 Stradivarius")

$newNames = @("23 - 9 Mango", "23 - 9 Zara", "23 - 9 Stradivarius", "23 - 9 Bershka")

foreach ($name in $newNames) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
    $newSheet.Name = $name
    $template.Range("A1:D1").Copy($newSheet.Range("A1:D1"))
}
